# Update column F (dSF) values for the rows that were re-pulled / recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    4  = 3
    6  = 0
    9  = 0
    10 = 2
    12 = 3
    14 = -7
    15 = 0
    21 = 1
    29 = -3
    37 = 1
    39 = 2
    41 = 3
    45 = 3
    53 = 3
    60 = -1
    62 = -1
    64 = 2
    66 = -3
    70 = 2
    71 = -1
    72 = -1
    75 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
